$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before column B to make room for the new weekly date
# columns (Sep_08 .. Jun_16), pushing the existing Jun_09..Mar_10 columns
# (and all data under them) from B:V to K:AE.
$ws.Range("B1:J1").EntireColumn.Insert()

# Populate the new header cells with the new week-ending dates. They are
# entered from the newest-adjacent-to-existing-data column (J, right next to
# the old Jun_09 column) outward to the most recent date (B), matching how
# the shared-string table ends up ordered chronologically.
$ws.Range("J1").Value = "Jun_16"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("B1").Value = "Sep_08"

# Fill the new columns for every existing data row with the default "UN"
# rating value, matching each row's original extent.
$lastRow = 33

for ($r = 2; $r -le 29; $r++) {
    $ws.Range("B" + $r + ":J" + $r).Value = "UN"
}
$ws.Range("B30:J30").Value = "UN"
$ws.Range("B31:J31").Value = "UN"
$ws.Range("B32:J32").Value = "UN"
$ws.Range("B33:J33").Value = "UN"
